# Update the workbook: add "Binary data" Load-Time and Size tables below the
# existing "Text data" tables, matching the commit "updated with binary times".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Retitle the two existing ("Text data") section headers.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "Load Time from Disk Serialization – Text  data (secs)"
$ws.Range("A9").Value2 = "Size of Serialized Data File – Text data (GB)"

# ---------------------------------------------------------------------------
# 2. Clone the layout/formatting of the existing second block (rows 9-15,
#    the "Size of Serialized Data File" table) down onto the two new blocks
#    (rows 17-23 and 25-31) so borders/fonts/shading line up. Formats only -
#    values are overwritten explicitly below.
#
#    The title rows (17 and 25) are merged *before* the format copy so that
#    the engine doesn't split their borders into "merged cell" left/mid/
#    right variants - pasting formats onto an already-merged range keeps a
#    single uniform style across it, just like rows 1 and 9 already have.
# ---------------------------------------------------------------------------
$ws.Range("A17:D17").MergeCells = $true
$ws.Range("A25:D25").MergeCells = $true

$ws.Range("A9:D15").Copy() | Out-Null
$ws.Range("A17:D23").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:D15").Copy() | Out-Null
$ws.Range("A25:D31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. New block #1: "Load Time from Disk Serialization – Binary data (secs)"
# ---------------------------------------------------------------------------
$ws.Range("A17").Value2 = "Load Time from Disk Serialization – Binary data (secs)"

$ws.Range("B18").Value2 = "Mesh"
$ws.Range("C18").Value2 = "Mesh + Supp."
$ws.Range("D18").Value2 = "SemmedDB"

$ws.Range("A19").Value2 = "UA"
$ws.Range("B19").Value2 = 152.525
$ws.Range("C19").Value2 = 540.335
$ws.Range("D19").Value2 = 74.6807

$ws.Range("A20").Value2 = "BCA"
$ws.Range("B20").Value2 = 106.683
$ws.Range("C20").Value2 = 333.878
$ws.Range("D20").Value2 = 66.2949

$ws.Range("A21").Value2 = "BB"
$ws.Range("B21").Value2 = "NA"
$ws.Range("C21").Value2 = "NA"
$ws.Range("D21").Value2 = 59.7214

$ws.Range("A22").Value2 = "Huffman"
$ws.Range("B22").Value2 = 88.2539
$ws.Range("C22").Value2 = 289.556
$ws.Range("D22").Value2 = 79.3396

$ws.Range("A23").Value2 = "Optimal"
$ws.Range("B23").Value2 = 70.4635
$ws.Range("C23").Value2 = 187.176
$ws.Range("D23").Value2 = 58.7767

# Row 24 is a blank spacer row (kept empty, just re-shaded below).

# ---------------------------------------------------------------------------
# 4. New block #2: "Size of Serialized Data File – Binary data (GB)"
# ---------------------------------------------------------------------------
$ws.Range("A25").Value2 = "Size of Serialized Data File – Binary data (GB)"

$ws.Range("B26").Value2 = "Mesh"
$ws.Range("C26").Value2 = "Mesh + Supp."
$ws.Range("D26").Value2 = "SemmedDB"

$ws.Range("A27").Value2 = "UA"
$ws.Range("B27").Value2 = 3.0423783156
$ws.Range("C27").Value2 = 9.982511392
$ws.Range("D27").Value2 = 1.932985173

$ws.Range("A28").Value2 = "BCA"
$ws.Range("B28").Value2 = 2.161175274
$ws.Range("C28").Value2 = 6.587253146
$ws.Range("D28").Value2 = 1.791065978

$ws.Range("A29").Value2 = "BB"
$ws.Range("B29").Value2 = "NA"
$ws.Range("C29").Value2 = "NA"
$ws.Range("D29").Value2 = 1.662723961

$ws.Range("A30").Value2 = "Huffman"
$ws.Range("B30").Value2 = 2.049469209
$ws.Range("C30").Value2 = 5.3057694
$ws.Range("D30").Value2 = 2.575099527

$ws.Range("A31").Value2 = "Optimal"
$ws.Range("B31").Value2 = 1.484849567
$ws.Range("C31").Value2 = 3.522409574
$ws.Range("D31").Value2 = 1.798093953

# ---------------------------------------------------------------------------
# 5. Re-colour the new blocks (new green-ish fill, replacing the pink one
#    that was copied over in step 2) across rows 17-31.
# ---------------------------------------------------------------------------
$newBlock = $ws.Range("A17:D31")
$newBlock.Interior.Color = 13434828
$newBlock.Interior.PatternColor = 16777164

# ---------------------------------------------------------------------------
# 6. Number formats: block #1's data values (rows 19-23) display with
#    "General" formatting (keeping their natural decimal places) while
#    block #2's data values (rows 27-31) keep the "0.000" formatting
#    (3 decimal places) used elsewhere in the sheet. The "NA" placeholder
#    cells in row 21 keep the "0.000"-formatted style.
# ---------------------------------------------------------------------------
$ws.Range("B19:D20").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("B22:D23").NumberFormat = "General"

# Row 24 spacer: no borders, just the new fill (already applied above).
$ws.Range("A24:D24").Borders.LineStyle = -4142

Write-Host "Binary data tables added."
